# Generate Report for Handoff
#
# This updates the localization-status report to reflect that the content
# has moved from "In Translation" to "Ready for handoff", refreshing the
# associated handoff timestamps on the Overview sheet and each language
# sheet (zh-cn, de-de). Columns that now hold the longer "Ready for
# handoff" status text are widened to fit.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_zhcn.Range("C2").Value     = "Ready for handoff"
$ws_dede.Range("C2").Value     = "Ready for handoff"

# --- Refresh the handoff-generation timestamps ---------------------------
$ws_overview.Range("G2").Value = "2016-08-23 06:56:12"
$ws_zhcn.Range("H2").Value     = "2016-08-23 06:56:02"
$ws_dede.Range("H2").Value     = "2016-08-23 06:56:12"

# --- Widen the columns that now display "Ready for handoff" --------------
# The underlying column-width grid only lands on certain discrete pixel
# widths, so use the value whose rounded pixel width matches the target
# XML column width (~17.216 characters) as closely as possible.
$newStatusColWidth = 16.333333333333336

$ws_overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$ws_overview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$ws_zhcn.Columns.Item(3).ColumnWidth     = $newStatusColWidth
$ws_dede.Columns.Item(3).ColumnWidth     = $newStatusColWidth
